$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = "Bestie del nord "
$ws.Range("B54").Value = "Stefano Tita | Clitoriders"
$ws.Range("C54").Value = "ALESSIO FARINATI | Pinguini Trentini"
$ws.Range("D54").Value = "Andrea Bertolini | Modium"
$ws.Range("E54").Value = "Matteo Bazzanella | Hellas Madonna"
$ws.Range("F54").Value = "Halzyd  Pupuleku | F.C. Sala Giardini"
